$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values re-pulled/recalculated per commit message
# "repull data, push all data, mean calculation"
$updates = @{
    2  = -6
    3  = -1
    5  = -9
    6  = 2
    8  = 6
    10 = -9
    11 = -1
    12 = -1
    13 = 3
    15 = -1
    16 = -3
    18 = 0
    20 = -4
    21 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
